# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and size condition codes were renumbered:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 and S25 are unchanged.)
#
# These codes appear embedded as substrings inside several text columns
# (Condition, Filename_Left, Filename_Right, Distance, Size), e.g.
# "Face15_D51_S25" -> "Face15_D55_S25", "Fixation_D64_l.png" -> "Fixation_D69_l.png".
# Walk every used cell and apply the substring substitutions to any
# string-valued cell, leaving numeric/boolean cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = $ws.UsedRange.Rows.Count
$cols = $ws.UsedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $nv = $v.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
